$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.642.62'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.285.09'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.19%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.55%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '113.67'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +18.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '268.17'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.56%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.33%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.616'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '47.68'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +6.45%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.18%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.58'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +9.87%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.41%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.55'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.629.26'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.850'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.48%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.287.02'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.523.01'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.62%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.52'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +5.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.23'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.53'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '233.05'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.99%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.53'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +5.15%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +13.26%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.41'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.72%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '43.43'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +7.56%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.49%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.44%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '176.35'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.80%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.69'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0926'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.80%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.49'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.34%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.79%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.70'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +7.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.109'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.94%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.91'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +18.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0353'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.98%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '75.39'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +16.96%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.67%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.85%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '13.17'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +8.84%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +5.99%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.32%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.91'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +13.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.76'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0999'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '101.28'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.10%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.80%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +6.69%  '
